$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1093
$ws.Range("I58").Value = 434.85715
$ws.Range("K58").Value = 1304.57145
$ws.Range("M58").Value = -1154.57145
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("H132").Value = 2210.4075
$ws.Range("I132").Value = 2267.24
$ws.Range("K132").Value = 6801.719999999999
$ws.Range("M132").Value = -4271.719999999999
$ws.Range("H138").Value = 5277.205
$ws.Range("J138").Value = 4849.1665
$ws.Range("L138").Value = 14547.4995
$ws.Range("N138").Value = -24827.4995
$ws.Range("H141").Value = 4736.875
$ws.Range("I141").Value = 3715
$ws.Range("K141").Value = 11145
$ws.Range("M141").Value = -5965
$ws.Range("N95").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8970.4
$ws.Range("I32").Value = 6194.871
$ws.Range("K32").Value = 6194.871
$ws.Range("M32").Value = -5907.871
$ws.Range("H46").Value = 15852.5
$ws.Range("J46").Value = 12626.2
$ws.Range("L46").Value = 12626.2
$ws.Range("N46").Value = -13264.2
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("H74").Value = 1503.0358
$ws.Range("I74").Value = 949.88
$ws.Range("K74").Value = 949.88
$ws.Range("M74").Value = -75.88
$ws.Range("H77").Value = 1503.0358
$ws.Range("I77").Value = 949.88
$ws.Range("K77").Value = 4749.4
$ws.Range("M77").Value = -381.3999999999996
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("H112").Value = 25000
$ws.Range("J112").Value = 25000
$ws.Range("L112").Value = 25000
$ws.Range("N112").Value = -27954
$ws.Range("H114").Value = 30398
$ws.Range("J114").Value = 30398
$ws.Range("L114").Value = 30398
$ws.Range("N114").Value = -39076
$ws.Range("H132").Value = 1909.5
$ws.Range("I132").Value = 1863.2858
$ws.Range("J132").Value = 2125.1667
$ws.Range("K132").Value = 5589.857400000001
$ws.Range("L132").Value = 6375.500100000001
$ws.Range("M132").Value = -3059.857400000001
$ws.Range("N132").Value = -11435.5001
$ws.Range("M60").ClearContents()
$ws.Range("N109").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 29999
$ws.Range("I141").Value = 10000
$ws.Range("K141").Value = 10000
$ws.Range("M141").Value = -4820

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 181.61765
$ws.Range("I7").Value = 135
$ws.Range("K7").Value = 135
$ws.Range("M7").Value = -22
$ws.Range("H58").Value = 3102.5
$ws.Range("I58").Value = 1338.4445
$ws.Range("J58").Value = 4036.4119
$ws.Range("K58").Value = 1338.4445
$ws.Range("L58").Value = 4036.4119
$ws.Range("M58").Value = -1135.4445
$ws.Range("N58").Value = -4442.4119
$ws.Range("H132").Value = 2490.0303
$ws.Range("I132").Value = 2255.375
$ws.Range("K132").Value = 6766.125
$ws.Range("M132").Value = -4236.125
$ws.Range("H134").Value = 3299.0908
$ws.Range("I134").Value = 1942.2
$ws.Range("J134").Value = 4429.8335
$ws.Range("K134").Value = 5826.6
$ws.Range("L134").Value = 13289.5005
$ws.Range("M134").Value = -3291.6
$ws.Range("N134").Value = -18359.5005
$ws.Range("H136").Value = 3102.5
$ws.Range("I136").Value = 1338.4445
$ws.Range("J136").Value = 4036.4119
$ws.Range("K136").Value = 4015.3335
$ws.Range("L136").Value = 12109.2357
$ws.Range("M136").Value = -1465.3335
$ws.Range("N136").Value = -17209.2357

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 29049992
$ws.Range("I4").Value = 32467318
$ws.Range("J4").Value = 2704
$ws.Range("K4").Value = 97401954
$ws.Range("L4").Value = 8112
$ws.Range("M4").Value = -97401842
$ws.Range("N4").Value = -8336
$ws.Range("H113").Value = 1011.8571
$ws.Range("J113").Value = 1276.6
$ws.Range("L113").Value = 3829.8
$ws.Range("N113").Value = -8169.799999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1020143.44
$ws.Range("I11").Value = 1270250
$ws.Range("K11").Value = 1270250
$ws.Range("M11").Value = -1270111
$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31262
$ws.Range("H70").Value = 7308.0835
$ws.Range("I70").Value = 7274.5
$ws.Range("J70").Value = 7341.6665
$ws.Range("K70").Value = 7274.5
$ws.Range("L70").Value = 7341.6665
$ws.Range("M70").Value = -7004.5
$ws.Range("N70").Value = -7881.6665
$ws.Range("H73").Value = 7308.0835
$ws.Range("I73").Value = 7274.5
$ws.Range("J73").Value = 7341.6665
$ws.Range("K73").Value = 7274.5
$ws.Range("L73").Value = 7341.6665
$ws.Range("M73").Value = -6338.5
$ws.Range("N73").Value = -9213.666499999999
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("H123").Value = 53750.25
$ws.Range("J123").Value = 53750.25
$ws.Range("L123").Value = 53750.25
$ws.Range("N123").Value = -58650.25
$ws.Range("H126").Value = 4660.857
$ws.Range("I126").Value = 3306
$ws.Range("J126").Value = 5202.8
$ws.Range("K126").Value = 9918
$ws.Range("L126").Value = 15608.4
$ws.Range("M126").Value = -7448
$ws.Range("N126").Value = -20548.4
$ws.Range("H132").Value = 3460.2666
$ws.Range("I132").Value = 2457.7778
$ws.Range("K132").Value = 7373.3334
$ws.Range("M132").Value = -4843.3334
$ws.Range("N103").ClearContents()
$ws.Range("N111").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3924.6667
$ws.Range("I61").Value = 3924.6667
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3924.6667
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3722.6667
$ws.Range("H93").Value = 2011.3846
$ws.Range("I93").Value = 1894.909
$ws.Range("J93").Value = 2652
$ws.Range("K93").Value = 1894.909
$ws.Range("L93").Value = 2652
$ws.Range("M93").Value = -646.9090000000001
$ws.Range("N93").Value = -5148
$ws.Range("H113").Value = 3924.6667
$ws.Range("I113").Value = 3924.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3924.6667
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1754.6667
$ws.Range("H132").Value = 3632.2368
$ws.Range("I132").Value = 2827.3914
$ws.Range("J132").Value = 4866.3335
$ws.Range("K132").Value = 8482.174199999999
$ws.Range("L132").Value = 14599.0005
$ws.Range("M132").Value = -5952.174199999999
$ws.Range("N132").Value = -19659.0005
$ws.Range("H136").Value = 3277.8
$ws.Range("J136").Value = 3699
$ws.Range("L136").Value = 11097
$ws.Range("N136").Value = -16197
$ws.Range("N61").ClearContents()
$ws.Range("N113").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5771.75
$ws.Range("J81").Value = 790
$ws.Range("L81").Value = 1580
$ws.Range("N81").Value = -3702
$ws.Range("H84").Value = 5771.75
$ws.Range("J84").Value = 790
$ws.Range("L84").Value = 7900
$ws.Range("N84").Value = -18508
$ws.Range("H132").Value = 1583.0667
$ws.Range("I132").Value = 1553.2858
$ws.Range("K132").Value = 4659.857400000001
$ws.Range("M132").Value = -2129.857400000001
$ws.Range("H136").Value = 3430.6943
$ws.Range("I136").Value = 1509.72
$ws.Range("K136").Value = 4529.16
$ws.Range("M136").Value = -1979.16

Write-Host "Applied 189 cell updates and 7 clears."